# Add a "Save" column (H) to the s_vals sheet, mirroring the header style
# used by the existing stat columns and filling in the save totals for
# each row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting of the neighboring header cell (G1) onto the new
# header cell (H1) so it picks up the same bold/centered/bordered style,
# then set its text.
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$ws.Range("H1").Value = "Save"

# Per-row "Save" values (column H), row 2 through row 21.
$saveValues = @(1, 1, 1, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 1, 0, 0, 1)

for ($i = 0; $i -lt $saveValues.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 8).Value = $saveValues[$i]
}
